$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.63867
$ws.Range("H2").Value = 1.91601
$ws.Range("I2").Value = 0.02162993170452444
$ws.Range("J2").Value = 0.02162993170452444
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 1.79832824002
$ws.Range("R2").Value = 16.18495416018
$ws.Range("S2").Value = 0.001665591037025534
$ws.Range("T2").Value = 0.001665591037025534

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.63867
$ws.Range("H3").Value = 1.91601
$ws.Range("I3").Value = 0.02162993170452444
$ws.Range("J3").Value = 0.02162993170452444
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("Q3").Value = 3.09281972287
$ws.Range("R3").Value = 27.83537750583
$ws.Range("S3").Value = 0.002864534235135393
$ws.Range("T3").Value = 0.002864534235135393

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.63867
$ws.Range("H4").Value = 1.91601
$ws.Range("I4").Value = 0.02162993170452444
$ws.Range("J4").Value = 0.02162993170452444
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 3.536820977969999
$ws.Range("R4").Value = 31.83138880173
$ws.Range("S4").Value = 0.003275763116751811
$ws.Range("T4").Value = 0.003275763116751812

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.63867
$ws.Range("H5").Value = 1.91601
$ws.Range("I5").Value = 0.02162993170452444
$ws.Range("J5").Value = 0.02162993170452444
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 0.75912039443
$ws.Range("R5").Value = 6.83208354987
$ws.Range("S5").Value = 0.0007030886224484994
$ws.Range("T5").Value = 0.0007030886224484995

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.63867
$ws.Range("H6").Value = 1.91601
$ws.Range("I6").Value = 0.02162993170452444
$ws.Range("J6").Value = 0.02162993170452444
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 12.03963182411
$ws.Range("R6").Value = 108.35668641699
$ws.Range("S6").Value = 0.01115096922189354
$ws.Range("T6").Value = 0.01115096922189354

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.63867
$ws.Range("H7").Value = 1.91601
$ws.Range("I7").Value = 0.02162993170452444
$ws.Range("J7").Value = 0.02162993170452444
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 2.12698100954
$ws.Range("R7").Value = 19.14282908586
$ws.Range("S7").Value = 0.001969985471269664
$ws.Range("T7").Value = 0.001969985471269664

# Row 8
$ws.Range("I8").Value = 0.9490095874171892
$ws.Range("J8").Value = 0.9490095874171893
$ws.Range("M8").Value = 2.815739333333333
$ws.Range("N8").Value = 8.447217999999999
$ws.Range("O8").Value = 0.07700398964630729
$ws.Range("P8").Value = 0.07700398964630729
$ws.Range("Q8").Value = 78.90134672709465
$ws.Range("R8").Value = 710.1121205438519
$ws.Range("S8").Value = 0.07307752444371959
$ws.Range("T8").Value = 0.07307752444371961

# Row 9
$ws.Range("I9").Value = 0.9490095874171892
$ws.Range("J9").Value = 0.9490095874171893
$ws.Range("O9").Value = 0.1324338085883186
$ws.Range("P9").Value = 0.1324338085883186
$ws.Range("S9").Value = 0.1256809540484872
$ws.Range("T9").Value = 0.1256809540484872

# Row 10
$ws.Range("I10").Value = 0.9490095874171892
$ws.Range("J10").Value = 0.9490095874171893
$ws.Range("M10").Value = 5.537790999999999
$ws.Range("N10").Value = 16.613373
$ws.Range("O10").Value = 0.1514458372546134
$ws.Range("P10").Value = 0.1514458372546134
$ws.Range("Q10").Value = 155.177420942558
$ws.Range("R10").Value = 1396.596788483022
$ws.Range("S10").Value = 0.1437235515290515
$ws.Range("T10").Value = 0.1437235515290515

# Row 11
$ws.Range("I11").Value = 0.9490095874171892
$ws.Range("J11").Value = 0.9490095874171893
$ws.Range("M11").Value = 1.188595666666667
$ws.Range("N11").Value = 3.565787
$ws.Range("O11").Value = 0.03250535563648733
$ws.Range("P11").Value = 0.03250535563648733
$ws.Range("Q11").Value = 33.30627864013533
$ws.Range("R11").Value = 299.756507761218
$ws.Range("S11").Value = 0.03084789414143185
$ws.Range("T11").Value = 0.03084789414143185

# Row 12
$ws.Range("I12").Value = 0.9490095874171892
$ws.Range("J12").Value = 0.9490095874171893
$ws.Range("M12").Value = 18.85109966666667
$ws.Range("N12").Value = 56.553299
$ws.Range("O12").Value = 0.5155341854158992
$ws.Range("P12").Value = 0.5155341854158992
$ws.Range("Q12").Value = 528.2368056512872
$ws.Range("R12").Value = 4754.131250861586
$ws.Range("S12").Value = 0.4892468846009992
$ws.Range("T12").Value = 0.4892468846009993

# Row 13
$ws.Range("I13").Value = 0.9490095874171892
$ws.Range("J13").Value = 0.9490095874171893
$ws.Range("M13").Value = 3.330328666666666
$ws.Range("N13").Value = 9.990985999999999
$ws.Range("O13").Value = 0.09107682345837424
$ws.Range("P13").Value = 0.09107682345837424
$ws.Range("Q13").Value = 93.32093128548931
$ws.Range("R13").Value = 839.8883815694039
$ws.Range("S13").Value = 0.08643277865349992
$ws.Range("T13").Value = 0.08643277865349994

# Row 14
$ws.Range("G14").Value = 0.866931
$ws.Range("H14").Value = 2.600793
$ws.Range("I14").Value = 0.02936048087828625
$ws.Range("J14").Value = 0.02936048087828625
$ws.Range("M14").Value = 2.815739333333333
$ws.Range("N14").Value = 8.447217999999999
$ws.Range("O14").Value = 0.07700398964630729
$ws.Range("P14").Value = 0.07700398964630729
$ws.Range("Q14").Value = 2.441051715986
$ws.Range("R14").Value = 21.969465443874
$ws.Range("S14").Value = 0.002260874165562157
$ws.Range("T14").Value = 0.002260874165562157

# Row 15
$ws.Range("G15").Value = 0.866931
$ws.Range("H15").Value = 2.600793
$ws.Range("I15").Value = 0.02936048087828625
$ws.Range("J15").Value = 0.02936048087828625
$ws.Range("O15").Value = 0.1324338085883186
$ws.Range("P15").Value = 0.1324338085883186
$ws.Range("Q15").Value = 4.198195147991
$ws.Range("R15").Value = 37.78375633191899
$ws.Range("S15").Value = 0.003888320304695948
$ws.Range("T15").Value = 0.003888320304695948

# Row 16
$ws.Range("G16").Value = 0.866931
$ws.Range("H16").Value = 2.600793
$ws.Range("I16").Value = 0.02936048087828625
$ws.Range("J16").Value = 0.02936048087828625
$ws.Range("M16").Value = 5.537790999999999
$ws.Range("N16").Value = 16.613373
$ws.Range("O16").Value = 0.1514458372546134
$ws.Range("P16").Value = 0.1514458372546134
$ws.Range("Q16").Value = 4.800882689421
$ws.Range("R16").Value = 43.20794420478899
$ws.Range("S16").Value = 0.004446522608810129
$ws.Range("T16").Value = 0.004446522608810129

# Row 17
$ws.Range("G17").Value = 0.866931
$ws.Range("H17").Value = 2.600793
$ws.Range("I17").Value = 0.02936048087828625
$ws.Range("J17").Value = 0.02936048087828625
$ws.Range("M17").Value = 1.188595666666667
$ws.Range("N17").Value = 3.565787
$ws.Range("O17").Value = 0.03250535563648733
$ws.Range("P17").Value = 0.03250535563648733
$ws.Range("Q17").Value = 1.030430429899
$ws.Range("R17").Value = 9.273873869091
$ws.Range("S17").Value = 0.0009543728726069803
$ws.Range("T17").Value = 0.0009543728726069803

# Row 18
$ws.Range("G18").Value = 0.866931
$ws.Range("H18").Value = 2.600793
$ws.Range("I18").Value = 0.02936048087828625
$ws.Range("J18").Value = 0.02936048087828625
$ws.Range("M18").Value = 18.85109966666667
$ws.Range("N18").Value = 56.553299
$ws.Range("O18").Value = 0.5155341854158992
$ws.Range("P18").Value = 0.5155341854158992
$ws.Range("Q18").Value = 16.342602685123
$ws.Range("R18").Value = 147.083424166107
$ws.Range("S18").Value = 0.01513633159300639
$ws.Range("T18").Value = 0.01513633159300639

# Row 19
$ws.Range("G19").Value = 0.866931
$ws.Range("H19").Value = 2.600793
$ws.Range("I19").Value = 0.02936048087828625
$ws.Range("J19").Value = 0.02936048087828625
$ws.Range("M19").Value = 3.330328666666666
$ws.Range("N19").Value = 9.990985999999999
$ws.Range("O19").Value = 0.09107682345837424
$ws.Range("P19").Value = 0.09107682345837424
$ws.Range("Q19").Value = 2.887165161322
$ws.Range("R19").Value = 25.984486451898
$ws.Range("S19").Value = 0.002674059333604649
$ws.Range("T19").Value = 0.002674059333604649
